$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.219.65"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "1.439.08"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'0.9069"
$ws.Range("E5").Value = "  -9.79%  "
$ws.Range("D6").Value = "'277.34"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D7").Value = "'0.3648"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3110"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").Value = "'39.07"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +5.27%  "
$ws.Range("D11").Value = "'0.06513"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'5.376"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "'17.55"
$ws.Range("E14").Value = "  +7.19%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "1.438.12"
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").Value = "'0.9437"
$ws.Range("E18").Value = "  -6.12%  "
$ws.Range("D19").Value = "'0.05632"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'67.66"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'5.386"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").Value = "'10.77"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").Value = "'2.267"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "20.266.05"
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "'2.158"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "'138.07"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "'16.90"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "1.590.50"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").Value = "'109.74"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "'3.814"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'0.8000"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").Value = "'4.808"
$ws.Range("E33").Value = "  -7.97%  "
$ws.Range("D34").Value = "'0.07679"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").Value = "'0.05919"
$ws.Range("E35").Value = "  +6.90%  "
$ws.Range("D36").Value = "'1.445"
$ws.Range("E36").Value = "  +12.21%  "
$ws.Range("D37").Value = "'4.654"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'1.130"
$ws.Range("E38").Value = "  +9.96%  "
$ws.Range("D39").Value = "'0.01980"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "'10.12"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'0.1833"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "'0.9203"
$ws.Range("E42").Value = "  -8.35%  "
$ws.Range("D43").Value = "'7.093"
$ws.Range("E43").Value = "  -13.16%  "
$ws.Range("D44").Value = "'3.514"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").Value = "'11.95"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("D47").Value = "'119.21"
$ws.Range("E47").Value = "  +10.45%  "
$ws.Range("D48").Value = "'0.5118"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").Value = "'1.754"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").Value = "'0.06327"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("D51").Value = "'0.9899"
$ws.Range("E51").Value = "  -1.62%  "

$ws.Range("D4,D5,D6,D7,D8,D9,D11,D12,D13,D14,D18,D19,D20,D22,D23,D24,D26,D27,D28,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D46,D47,D48,D49,D50,D51").Style = "Normal"
